$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TABLE_1")
$ws2 = $wb.Worksheets.Item("TABLE_2")

# --- Header row 4: new date columns (force text format so dates stay as strings) ---
$ws1.Range("EH4:EJ4").NumberFormat = "@"
$ws1.Range("EH4").Value = "05/01/2023"
$ws1.Range("EI4").Value = "06/01/2023"
$ws1.Range("EJ4").Value = "07/01/2023"
$ws2.Range("DV4:DX4").NumberFormat = "@"
$ws2.Range("DV4").Value = "05/01/2023"
$ws2.Range("DW4").Value = "06/01/2023"
$ws2.Range("DX4").Value = "07/01/2023"

# --- TABLE_1 (sheet1) data rows ---
$ws1.Range("EF5").Value = 10817.6
$ws1.Range("EG5").Value = 10826.6
$ws1.Range("EH5").Value = 10720.3
$ws1.Range("EI5").Value = 10139.8
$ws1.Range("EJ5").Value = 9046.6
$ws1.Range("EG6").Value = 171.6
$ws1.Range("EH6").Value = 171.5
$ws1.Range("EI6").Value = 169.2
$ws1.Range("EJ6").Value = 160.8
$ws1.Range("EH7").Value = 28.3
$ws1.Range("EI7").Value = 22.9
$ws1.Range("EJ7").Value = 19.3
$ws1.Range("EG8").Value = 199.4
$ws1.Range("EH8").Value = 193.5
$ws1.Range("EI8").Value = 152.7
$ws1.Range("EJ8").Value = 151.9
$ws1.Range("EH9").Value = 101.6
$ws1.Range("EI9").Value = 93.4
$ws1.Range("EJ9").Value = 82.2
$ws1.Range("EH10").Value = 1255.9
$ws1.Range("EI10").Value = 1234.6
$ws1.Range("EJ10").Value = 1074.5
$ws1.Range("EH11").Value = 225.6
$ws1.Range("EI11").Value = 211.2
$ws1.Range("EJ11").Value = 197
$ws1.Range("EG12").Value = 124.8
$ws1.Range("EH12").Value = 123
$ws1.Range("EI12").Value = 114.2
$ws1.Range("EJ12").Value = 96.7
$ws1.Range("EG13").Value = 38.2
$ws1.Range("EH13").Value = 37.6
$ws1.Range("EI13").Value = 35.2
$ws1.Range("EJ13").Value = 34.1
$ws1.Range("EG15").Value = 488.9
$ws1.Range("EH15").Value = 480.4
$ws1.Range("EI15").Value = 399.2
$ws1.Range("EJ15").Value = 393
$ws1.Range("EG16").Value = 346.3
$ws1.Range("EH16").Value = 342.7
$ws1.Range("EI16").Value = 335.3
$ws1.Range("EJ16").Value = 320
$ws1.Range("EH18").Value = 64.5
$ws1.Range("EI18").Value = 61.4
$ws1.Range("EJ18").Value = 54.9
$ws1.Range("EG19").Value = 446
$ws1.Range("EH19").Value = 445
$ws1.Range("EI19").Value = 408.7
$ws1.Range("EJ19").Value = 380.6
$ws1.Range("EG20").Value = 229
$ws1.Range("EH20").Value = 229.5
$ws1.Range("EI20").Value = 185.5
$ws1.Range("EJ20").Value = 167.8
$ws1.Range("EG21").Value = 146.2
$ws1.Range("EH21").Value = 145.4
$ws1.Range("EI21").Value = 132.2
$ws1.Range("EJ21").Value = 114.6
$ws1.Range("EG22").Value = 135.6
$ws1.Range("EH22").Value = 135.3
$ws1.Range("EI22").Value = 119.9
$ws1.Range("EJ22").Value = 102.6
$ws1.Range("EG23").Value = 161
$ws1.Range("EH23").Value = 160.7
$ws1.Range("EI23").Value = 151.1
$ws1.Range("EJ23").Value = 126.9
$ws1.Range("EG24").Value = 150.3
$ws1.Range("EH24").Value = 149
$ws1.Range("EI24").Value = 139.1
$ws1.Range("EJ24").Value = 131.9
$ws1.Range("EG25").Value = 49.4
$ws1.Range("EH25").Value = 49.7
$ws1.Range("EI25").Value = 47.7
$ws1.Range("EJ25").Value = 40.7
$ws1.Range("EH26").Value = 225.2
$ws1.Range("EI26").Value = 209.7
$ws1.Range("EJ26").Value = 196.1
$ws1.Range("EG27").Value = 245.9
$ws1.Range("EH27").Value = 245.1
$ws1.Range("EI27").Value = 239.7
$ws1.Range("EJ27").Value = 212.6
$ws1.Range("EG28").Value = 324.6
$ws1.Range("EH28").Value = 308.1
$ws1.Range("EI28").Value = 291.9
$ws1.Range("EJ28").Value = 263.6
$ws1.Range("EG29").Value = 205.9
$ws1.Range("EH29").Value = 203
$ws1.Range("EI29").Value = 197.7
$ws1.Range("EJ29").Value = 167.5
$ws1.Range("EH30").Value = 101.9
$ws1.Range("EI30").Value = 98.7
$ws1.Range("EJ30").Value = 93.3
$ws1.Range("EG32").Value = 43.6
$ws1.Range("EH32").Value = 42.7
$ws1.Range("EI32").Value = 39.1
$ws1.Range("EJ32").Value = 32.1
$ws1.Range("EH33").Value = 94.7
$ws1.Range("EI33").Value = 90.5
$ws1.Range("EJ33").Value = 78.7
$ws1.Range("EH34").Value = 83.9
$ws1.Range("EI34").Value = 76.5
$ws1.Range("EJ34").Value = 70.6
$ws1.Range("EG35").Value = 46.4
$ws1.Range("EH35").Value = 45.5
$ws1.Range("EI35").Value = 39.1
$ws1.Range("EJ35").Value = 32.3
$ws1.Range("EG36").Value = 312.6
$ws1.Range("EH36").Value = 310
$ws1.Range("EI36").Value = 303.6
$ws1.Range("EJ36").Value = 257.1
$ws1.Range("EH37").Value = 77.1
$ws1.Range("EI37").Value = 71
$ws1.Range("EJ37").Value = 65.3
$ws1.Range("EG38").Value = 676.1
$ws1.Range("EH38").Value = 672.4
$ws1.Range("EI38").Value = 659.8
$ws1.Range("EJ38").Value = 538.4
$ws1.Range("EH39").Value = 333.1
$ws1.Range("EI39").Value = 293.8
$ws1.Range("EJ39").Value = 237.3
$ws1.Range("EG40").Value = 38.2
$ws1.Range("EH40").Value = 37.9
$ws1.Range("EI40").Value = 32.4
$ws1.Range("EJ40").Value = 26.1
$ws1.Range("EH41").Value = 385.3
$ws1.Range("EI41").Value = 357.7
$ws1.Range("EJ41").Value = 348.9
$ws1.Range("EG42").Value = 147.7
$ws1.Range("EH42").Value = 146.8
$ws1.Range("EI42").Value = 134.8
$ws1.Range("EJ42").Value = 125.6
$ws1.Range("EH43").Value = 142.2
$ws1.Range("EI43").Value = 141.6
$ws1.Range("EJ43").Value = 111.3
$ws1.Range("EG44").Value = 332
$ws1.Range("EH44").Value = 326.7
$ws1.Range("EI44").Value = 302.2
$ws1.Range("EJ44").Value = 271.3
$ws1.Range("EG45").Value = 29.5
$ws1.Range("EH45").Value = 29.8
$ws1.Range("EI45").Value = 27.3
$ws1.Range("EJ45").Value = 22.9
$ws1.Range("EH46").Value = 160.9
$ws1.Range("EI46").Value = 155.8
$ws1.Range("EJ46").Value = 147.1
$ws1.Range("EH47").Value = 37.7
$ws1.Range("EI47").Value = 34
$ws1.Range("EJ47").Value = 28.4
$ws1.Range("EG48").Value = 209.4
$ws1.Range("EH48").Value = 205.8
$ws1.Range("EI48").Value = 192.8
$ws1.Range("EJ48").Value = 168.2
$ws1.Range("EH49").Value = 1178.1
$ws1.Range("EI49").Value = 1133.6
$ws1.Range("EJ49").Value = 1061
$ws1.Range("EH50").Value = 128.6
$ws1.Range("EI50").Value = 119.4
$ws1.Range("EJ50").Value = 107.5
$ws1.Range("EH51").Value = 29.6
$ws1.Range("EI51").Value = 29.4
$ws1.Range("EJ51").Value = 26.1
$ws1.Range("EG52").Value = 323.7
$ws1.Range("EH52").Value = 313.6
$ws1.Range("EI52").Value = 305.8
$ws1.Range("EJ52").Value = 273.2
$ws1.Range("EG53").Value = 252.1
$ws1.Range("EH53").Value = 252.7
$ws1.Range("EI53").Value = 249.5
$ws1.Range("EJ53").Value = 230.2
$ws1.Range("EH54").Value = 61.6
$ws1.Range("EI54").Value = 55.9
$ws1.Range("EJ54").Value = 51.7
$ws1.Range("EH55").Value = 211.4
$ws1.Range("EI55").Value = 189.7
$ws1.Range("EJ55").Value = 168.1
$ws1.Range("EG56").Value = 30.6
$ws1.Range("EH56").Value = 30.7
$ws1.Range("EI56").Value = 28.1
$ws1.Range("EJ56").Value = 24.1

# --- TABLE_2 (sheet2) data rows ---
$ws2.Range("DT5").Value = 2.36765902681833
$ws2.Range("DU5").Value = 2.26894884002117
$ws2.Range("DV5").Value = 2.36326483843862
$ws2.Range("DW5").Value = 2.67112191170514
$ws2.Range("DX5").Value = 1.785573645068
$ws2.Range("DU6").Value = 0.941176470588232
$ws2.Range("DV6").Value = 0.823045267489698
$ws2.Range("DW6").Value = 0.594530321046373
$ws2.Range("DX6").Value = 0.06222775357811
$ws2.Range("DV7").Value = 5.99250936329586
$ws2.Range("DW7").Value = 4.09090909090908
$ws2.Range("DX7").Value = 2.11640211640213
$ws2.Range("DU8").Value = 1.83861082737485
$ws2.Range("DV8").Value = 1.94942044257112
$ws2.Range("DW8").Value = 1.73217854763491
$ws2.Range("DX8").Value = 3.47411444141689
$ws2.Range("DV9").Value = 0.594059405940588
$ws2.Range("DW9").Value = 0.755124056094948
$ws2.Range("DX9").Value = 0.611995104039185
$ws2.Range("DV10").Value = 2.75732286041565
$ws2.Range("DW10").Value = 2.72069223729094
$ws2.Range("DX10").Value = 3.08932169241102
$ws2.Range("DV11").Value = 3.15500685871058
$ws2.Range("DW11").Value = 3.32681017612524
$ws2.Range("DX11").Value = 5.74342458400429
$ws2.Range("DU12").Value = 0.402252614641995
$ws2.Range("DV12").Value = 0.408163265306122
$ws2.Range("DW12").Value = 2.42152466367713
$ws2.Range("DX12").Value = -0.309278350515461
$ws2.Range("DU13").Value = 1.05820105820107
$ws2.Range("DV13").Value = 1.89701897018969
$ws2.Range("DW13").Value = 0.57142857142858
$ws2.Range("DX13").Value = 0.887573964497054
$ws2.Range("DU15").Value = 1.76935886761032
$ws2.Range("DV15").Value = 1.50010564124233
$ws2.Range("DW15").Value = 1.91473066122032
$ws2.Range("DX15").Value = 1.83985488468516
$ws2.Range("DU16").Value = 2.82066508313539
$ws2.Range("DV16").Value = 2.02441202738911
$ws2.Range("DW16").Value = 3.51960481630134
$ws2.Range("DX16").Value = 3.55987055016181
$ws2.Range("DV18").Value = 5.73770491803279
$ws2.Range("DW18").Value = 5.67986230636834
$ws2.Range("DX18").Value = 6.60194174757281
$ws2.Range("DU19").Value = 3.89005357558816
$ws2.Range("DV19").Value = 4.02057035998131
$ws2.Range("DW19").Value = 3.65204159269591
$ws2.Range("DX19").Value = 3.42391304347827
$ws2.Range("DU20").Value = 2.18652387327086
$ws2.Range("DV20").Value = 2.27272727272727
$ws2.Range("DW20").Value = -2.87958115183246
$ws2.Range("DX20").Value = -1.46799765120376
$ws2.Range("DU21").Value = 2.23776223776223
$ws2.Range("DV21").Value = 1.89208128941837
$ws2.Range("DW21").Value = 0.379650721336371
$ws2.Range("DX21").Value = 0
$ws2.Range("DU22").Value = -0.147275405007376
$ws2.Range("DV22").Value = -0.147601476014773
$ws2.Range("DW22").Value = -0.580431177446093
$ws2.Range("DX22").Value = 0.686947988223738
$ws2.Range("DU23").Value = 1.70562223626026
$ws2.Range("DV23").Value = 1.64452877925363
$ws2.Range("DW23").Value = 1.75084175084175
$ws2.Range("DX23").Value = 3.42298288508559
$ws2.Range("DU24").Value = 3.29896907216496
$ws2.Range("DV24").Value = 2.40549828178694
$ws2.Range("DW24").Value = 3.34323922734027
$ws2.Range("DX24").Value = 2.08978328173376
$ws2.Range("DU25").Value = -1.39720558882235
$ws2.Range("DV25").Value = 1.22199592668026
$ws2.Range("DW25").Value = 0.632911392405072
$ws2.Range("DX25").Value = 0.742574257425753
$ws2.Range("DV26").Value = 3.01921317474838
$ws2.Range("DW26").Value = 2.04379562043796
$ws2.Range("DX26").Value = 2.56276150627615
$ws2.Range("DU27").Value = 3.05951383067894
$ws2.Range("DV27").Value = 2.50941028858217
$ws2.Range("DW27").Value = 2.43589743589743
$ws2.Range("DX27").Value = 2.70531400966183
$ws2.Range("DU28").Value = 6.07843137254903
$ws2.Range("DV28").Value = 6.24137931034484
$ws2.Range("DW28").Value = 5.6842867487328
$ws2.Range("DX28").Value = 5.73606097071801
$ws2.Range("DU29").Value = 0.292255236239646
$ws2.Range("DV29").Value = 0.644521566683198
$ws2.Range("DW29").Value = 1.07361963190182
$ws2.Range("DX29").Value = -0.475341651812247
$ws2.Range("DV30").Value = 1.09126984126985
$ws2.Range("DW30").Value = 1.64778578784757
$ws2.Range("DX30").Value = -1.06044538706257
$ws2.Range("DU32").Value = -1.35746606334842
$ws2.Range("DV32").Value = -4.26008968609865
$ws2.Range("DW32").Value = -6.9047619047619
$ws2.Range("DX32").Value = -7.49279538904899
$ws2.Range("DV33").Value = 1.39186295503212
$ws2.Range("DW33").Value = 2.95790671217292
$ws2.Range("DX33").Value = 1.94300518134717
$ws2.Range("DV34").Value = 4.74406991260925
$ws2.Range("DW34").Value = 7.74647887323944
$ws2.Range("DX34").Value = 5.68862275449101
$ws2.Range("DU35").Value = 0.432900432900424
$ws2.Range("DV35").Value = 0.886917960088689
$ws2.Range("DW35").Value = 0.25641025641026
$ws2.Range("DX35").Value = -0.0000000000000219982271133158
$ws2.Range("DU36").Value = 0.806191551112544
$ws2.Range("DV36").Value = 1.57273918741809
$ws2.Range("DW36").Value = 0.729927007299285
$ws2.Range("DX36").Value = 3.66935483870969
$ws2.Range("DV37").Value = 3.35120643431635
$ws2.Range("DW37").Value = 11.6352201257862
$ws2.Range("DX37").Value = 7.75577557755775
$ws2.Range("DU38").Value = 3.52166590108712
$ws2.Range("DV38").Value = 3.49392027089425
$ws2.Range("DW38").Value = 3.33594361785436
$ws2.Range("DX38").Value = 1.7961807525052
$ws2.Range("DV39").Value = 0.725733293014827
$ws2.Range("DW39").Value = 0.789022298456265
$ws2.Range("DX39").Value = 0.721561969439736
$ws2.Range("DU40").Value = 1.3262599469496
$ws2.Range("DV40").Value = 0.530503978779829
$ws2.Range("DW40").Value = -2.7027027027027
$ws2.Range("DX40").Value = -6.4516129032258
$ws2.Range("DV41").Value = -0.874710573707223
$ws2.Range("DW41").Value = 0.562271577171774
$ws2.Range("DX41").Value = 1.39494333042719
$ws2.Range("DU42").Value = -1.07166778298729
$ws2.Range("DV42").Value = -1.21130551816957
$ws2.Range("DW42").Value = 1.73584905660378
$ws2.Range("DX42").Value = 2.61437908496731
$ws2.Range("DV43").Value = 5.17751479289943
$ws2.Range("DW43").Value = 7.19152157456472
$ws2.Range("DX43").Value = 2.39190432382705
$ws2.Range("DU44").Value = 1.96560196560196
$ws2.Range("DV44").Value = 2.02998126171143
$ws2.Range("DW44").Value = 1.07023411371237
$ws2.Range("DX44").Value = 0.930059523809524
$ws2.Range("DU45").Value = 1.02739726027398
$ws2.Range("DV45").Value = 2.75862068965516
$ws2.Range("DW45").Value = 1.48698884758365
$ws2.Range("DX45").Value = -2.13675213675214
$ws2.Range("DV46").Value = 0
$ws2.Range("DW46").Value = 0.257400257400261
$ws2.Range("DX46").Value = 0
$ws2.Range("DV47").Value = 1.8918918918919
$ws2.Range("DW47").Value = 2.10210210210211
$ws2.Range("DX47").Value = 1.06761565836298
$ws2.Range("DU48").Value = 1.35527589545014
$ws2.Range("DV48").Value = 2.74588117823265
$ws2.Range("DW48").Value = 0.626304801670155
$ws2.Range("DX48").Value = 0.899820035992801
$ws2.Range("DV49").Value = 3.07988450433107
$ws2.Range("DW49").Value = 2.31046931407941
$ws2.Range("DX49").Value = 1.50196115947575
$ws2.Range("DV50").Value = 3.04487179487179
$ws2.Range("DW50").Value = 3.10880829015544
$ws2.Range("DX50").Value = 6.01577909270216
$ws2.Range("DV51").Value = -5.12820512820513
$ws2.Range("DW51").Value = 8.08823529411764
$ws2.Range("DX51").Value = 5.24193548387099
$ws2.Range("DU52").Value = 2.92527821939586
$ws2.Range("DV52").Value = 2.01691607026677
$ws2.Range("DW52").Value = 2.27424749163878
$ws2.Range("DX52").Value = 3.09433962264151
$ws2.Range("DU53").Value = 2.23033252230333
$ws2.Range("DV53").Value = 2.10101010101011
$ws2.Range("DW53").Value = 3.35542667771334
$ws2.Range("DX53").Value = 0.43630017452007
$ws2.Range("DV54").Value = 1.14942528735631
$ws2.Range("DW54").Value = -0.356506238859172
$ws2.Range("DX54").Value = 1.57170923379176
$ws2.Range("DU55").Value = 2.88184438040346
$ws2.Range("DV55").Value = 3.62745098039216
$ws2.Range("DW55").Value = 1.77038626609441
$ws2.Range("DX55").Value = -1.05944673337256
$ws2.Range("DU56").Value = 2
$ws2.Range("DV56").Value = 2.33333333333333
$ws2.Range("DW56").Value = 2.93040293040293
$ws2.Range("DX56").Value = 0
